# Enable maps for dwellings, income and access.
#
# The "menu" sheet drives which visualizations show up in the SILO
# visualizer. This adds two new spatial "Accessibilities" map entries
# (Auto_accessibility / Transit_accessibility) just under the existing
# "Accessibilities" header row, shifting the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("menu")
$cfg = $wb.Worksheets.Item("config")

# Row 12 ("spatial" / "Accessibilities" / "accessibilities") currently has
# no category columns filled in -- give it the "Auto_accessibility" map.
$ws.Range("D12").Value() = "Auto_accessibility"

# Insert a fresh row right after it for the "Transit_accessibility" map;
# this pushes every row from the old 13 onward down by one.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value() = "spatial"
$ws.Range("B13").Value() = "Accessibilities"
$ws.Range("C13").Value() = "accessibilities"
$ws.Range("D13").Value() = "Transit_accessibility"
$ws.Range("G13").Value() = "yes"

# Fill in the category_value codes last so the new shared-string entries
# land in the same order as the source edit.
$ws.Range("E12").Value() = "autoAccessibility"
$ws.Range("E13").Value() = "transitAccessibility"

# Restore the selection/active-sheet state left behind by the edit.
$cfg.Activate()
$cfg.Range("C9").Select()

$ws.Activate()
$ws.Range("E14").Select()
